$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5041017266062511
$ws.Range("C2").Value = 0.1581249250453922
$ws.Range("E2").Value = 0.2165337921211758
$ws.Range("F2").Value = 2.231111514047953
$ws.Range("G2").Value = 0.002494084295733952
$ws.Range("I2").Value = 1.147359092487445
$ws.Range("J2").Value = 0.08197583239966022
$ws.Range("K2").Value = 0.2651977151940059
$ws.Range("L2").Value = 0.4315841589738909
$ws.Range("M2").Value = 0.2178704624623471
$ws.Range("O2").Value = 4.295835490669134

$ws.Range("B3").Value = 0.4689878571626593
$ws.Range("C3").Value = 0.157603835916138
$ws.Range("E3").Value = 0.2164078504903095
$ws.Range("F3").Value = 2.23631421078634
$ws.Range("G3").Value = 0.002496355630131978
$ws.Range("I3").Value = 1.157062018474946
$ws.Range("J3").Value = 0.08019814982303686
$ws.Range("K3").Value = 0.2347821990510397
$ws.Range("L3").Value = 0.4263310002189371
$ws.Range("M3").Value = 0.209770911978918
$ws.Range("O3").Value = 4.328891240760242

$ws.Range("B4").Value = 0.4475362412818811
$ws.Range("C4").Value = 0.1572862086821196
$ws.Range("E4").Value = 0.2163989054339162
$ws.Range("F4").Value = 2.240513565828707
$ws.Range("G4").Value = 0.002497824995069194
$ws.Range("I4").Value = 1.163497698915101
$ws.Range("J4").Value = 0.07909739515754666
$ws.Range("K4").Value = 0.2160919838259474
$ws.Range("L4").Value = 0.4232623875059716
$ws.Range("M4").Value = 0.2048696086078543
$ws.Range("O4").Value = 4.351138808109411

$ws.Range("B5").Value = 0.4388224589499714
$ws.Range("C5").Value = 0.1571573711001619
$ws.Range("E5").Value = 0.2164125085467354
$ws.Range("F5").Value = 2.242477846069917
$ws.Range("G5").Value = 0.002498442625987141
$ws.Range("I5").Value = 1.1662405683104
$ws.Range("J5").Value = 0.07864652464418853
$ws.Range("K5").Value = 0.2084722786108415
$ws.Range("L5").Value = 0.4220515205571473
$ws.Range("M5").Value = 0.2028905316052345
$ws.Range("O5").Value = 4.360695742804182

$ws.Range("B6").Value = 0.4373772480135187
$ws.Range("C6").Value = 0.1571360142498399
$ws.Range("E6").Value = 0.2164158106367822
$ws.Range("F6").Value = 2.242819303954199
$ws.Range("G6").Value = 0.002498546323315325
$ws.Range("I6").Value = 1.166703286176745
$ws.Range("J6").Value = 0.07857151942842222
$ws.Range("K6").Value = 0.2072068475479512
$ws.Range("L6").Value = 0.421852855037713
$ws.Range("M6").Value = 0.2025630144428376
$ws.Range("O6").Value = 4.362312318508742

$ws.Range("B7").Value = 0.4474186102255828
$ws.Range("C7").Value = 0.1572844686937849
$ws.Range("E7").Value = 0.2163990189824041
$ws.Range("F7").Value = 2.240539031945467
$ws.Range("G7").Value = 0.002497833248032367
$ws.Range("I7").Value = 1.163534203121241
$ws.Range("J7").Value = 0.07909132386001261
$ws.Range("K7").Value = 0.2159892345133301
$ws.Range("L7").Value = 0.4232458966851169
$ws.Range("M7").Value = 0.2048428439808205
$ws.Range("O7").Value = 4.351265708462265

$ws.Range("B8").Value = 0.4919723393412596
$ws.Range("C8").Value = 0.15794477931518
$ws.Range("E8").Value = 0.2164762056128744
$ws.Range("F8").Value = 2.232697032330165
$ws.Range("G8").Value = 0.002494851968934989
$ws.Range("I8").Value = 1.150605485008324
$ws.Range("J8").Value = 0.08136482009261314
$ws.Range("K8").Value = 0.2547138249566956
$ws.Range("L8").Value = 0.429740419718172
$ws.Range("M8").Value = 0.2150629267973727
$ws.Range("O8").Value = 4.306828275286122

$ws.Range("B9").Value = 0.5801782079079487
$ws.Range("C9").Value = 0.159257550996692
$ws.Range("E9").Value = 0.2171680172868911
$ws.Range("F9").Value = 2.225279608979648
$ws.Range("G9").Value = 0.002489596355754291
$ws.Range("I9").Value = 1.129042303357043
$ws.Range("J9").Value = 0.08574896509409058
$ws.Range("K9").Value = 0.330516835887039
$ws.Range("L9").Value = 0.4437138961751401
$ws.Range("M9").Value = 0.2356680298084513
$ws.Range("O9").Value = 4.235159321946867

$ws.Range("B10").Value = 0.6454668458731589
$ws.Range("C10").Value = 0.1602323537429129
$ws.Range("E10").Value = 0.2180031313145889
$ws.Range("F10").Value = 2.224668742281978
$ws.Range("G10").Value = 0.002486091642042253
$ws.Range("I10").Value = 1.115506131465903
$ws.Range("J10").Value = 0.08892401598714628
$ws.Range("K10").Value = 0.3861090661393405
$ws.Range("L10").Value = 0.4547271123236243
$ws.Range("M10").Value = 0.251143016762164
$ws.Range("O10").Value = 4.191925339261616

$ws.Range("B11").Value = 0.6752683115126672
$ws.Range("C11").Value = 0.1606779358271453
$ws.Range("E11").Value = 0.2184534750322094
$ws.Range("F11").Value = 2.225438649309964
$ws.Range("G11").Value = 0.002484573941860275
$ws.Range("I11").Value = 1.109848096546269
$ws.Range("J11").Value = 0.09035829889867131
$ws.Range("K11").Value = 0.4113742885546117
$ws.Range("L11").Value = 0.4598979595878916
$ws.Range("M11").Value = 0.2582546310596214
$ws.Range("O11").Value = 4.174300525314862

$ws.Range("B12").Value = 0.6865673496004661
$ws.Range("C12").Value = 0.1608469617876196
$ws.Range("E12").Value = 0.2186340888636913
$ws.Range("F12").Value = 2.225880581373431
$ws.Range("G12").Value = 0.002484010188618116
$ws.Range("I12").Value = 1.107777331429403
$ws.Range("J12").Value = 0.09089995831770636
$ws.Range("K12").Value = 0.4209377212879417
$ws.Range("L12").Value = 0.4618789999306188
$ws.Range("M12").Value = 0.2609578076485803
$ws.Range("O12").Value = 4.167920005564156

$ws.Range("B13").Value = 0.6841332938806772
$ws.Range("C13").Value = 0.160810546179718
$ws.Range("E13").Value = 0.2185947430239636
$ws.Range("F13").Value = 2.22577871955491
$ws.Range("G13").Value = 0.002484131116136734
$ws.Range("I13").Value = 1.108220114575268
$ws.Range("J13").Value = 0.0907833681580712
$ws.Range("K13").Value = 0.4188782491517316
$ws.Range("L13").Value = 0.4614513295120446
$ws.Range("M13").Value = 0.2603751811771886
$ws.Range("O13").Value = 4.169281107651074

$ws.Range("B14").Value = 0.6761976162195822
$ws.Range("C14").Value = 0.1606918358868654
$ws.Range("E14").Value = 0.2184681325494751
$ws.Range("F14").Value = 2.225471995031626
$ws.Range("G14").Value = 0.002484527341977162
$ws.Range("I14").Value = 1.109676294489354
$ws.Range("J14").Value = 0.09040289119915457
$ws.Range("K14").Value = 0.4121611601920563
$ws.Range("L14").Value = 0.4600604820027172
$ws.Range("M14").Value = 0.258476820414792
$ws.Range("O14").Value = 4.173769712210344

$ws.Range("B15").Value = 0.6713385699641776
$ws.Range("C15").Value = 0.1606191602416516
$ws.Range("E15").Value = 0.2183918909567417
$ws.Range("F15").Value = 2.225303693452233
$ws.Range("G15").Value = 0.002484771468467152
$ws.Range("I15").Value = 1.110577597398787
$ws.Range("J15").Value = 0.09016964557688567
$ws.Range("K15").Value = 0.4080462185174838
$ws.Range("L15").Value = 0.4592115314730592
$ws.Range("M15").Value = 0.2573153375911019
$ws.Range("O15").Value = 4.176557346450693

$ws.Range("B16").Value = 0.6435212290885772
$ws.Range("C16").Value = 0.1602032757632088
$ws.Range("E16").Value = 0.2179751134354788
$ws.Range("F16").Value = 2.224639486863808
$ws.Range("G16").Value = 0.002486192364824943
$ws.Range("I16").Value = 1.115885948139741
$ws.Range("J16").Value = 0.08883007756180916
$ws.Range("K16").Value = 0.3844573983765542
$ws.Range("L16").Value = 0.4543924070834464
$ws.Range("M16").Value = 0.250679688267617
$ws.Range("O16").Value = 4.193118254888674

$ws.Range("B17").Value = 0.6264816459799363
$ws.Range("C17").Value = 0.159948682085016
$ws.Range("E17").Value = 0.2177374401237167
$ws.Range("F17").Value = 2.224500179132207
$ws.Range("G17").Value = 0.002487083625066865
$ws.Range("I17").Value = 1.119270391882146
$ws.Range("J17").Value = 0.08800570047277745
$ws.Range("K17").Value = 0.3699799185682195
$ws.Range("L17").Value = 0.4514771063368386
$ws.Range("M17").Value = 0.2466272319721554
$ws.Range("O17").Value = 4.203800908320318

$ws.Range("B18").Value = 0.6166905086163013
$ws.Range("C18").Value = 0.1598024487357392
$ws.Range("E18").Value = 0.217607370531276
$ws.Range("F18").Value = 2.224518693456957
$ws.Range("G18").Value = 0.002487603468696403
$ws.Range("I18").Value = 1.12126406487398
$ws.Range("J18").Value = 0.08753059484675418
$ws.Range("K18").Value = 0.3616506302518303
$ws.Range("L18").Value = 0.4498154560995147
$ws.Range("M18").Value = 0.2443031510877205
$ws.Range("O18").Value = 4.210137547671962

$ws.Range("B19").Value = 0.6133770645689651
$ws.Range("C19").Value = 0.1597529717895512
$ws.Range("E19").Value = 0.2175644721392445
$ws.Range("F19").Value = 2.224541912395892
$ws.Range("G19").Value = 0.002487780719253626
$ws.Range("I19").Value = 1.121947166965327
$ws.Range("J19").Value = 0.08736957058419392
$ws.Range("K19").Value = 0.3588301098801026
$ws.Range("L19").Value = 0.4492554580678245
$ws.Range("M19").Value = 0.2435174281969452
$ws.Range("O19").Value = 4.212316048722471

$ws.Range("B20").Value = 0.6282945522738999
$ws.Range("C20").Value = 0.1599757631922429
$ws.Range("E20").Value = 0.2177620546534627
$ws.Range("F20").Value = 2.224504801833575
$ws.Range("G20").Value = 0.002486988002816427
$ws.Range("I20").Value = 1.118905244861779
$ws.Range("J20").Value = 0.0880935548831232
$ws.Range("K20").Value = 0.3715213049926547
$ws.Range("L20").Value = 0.4517858780490513
$ws.Range("M20").Value = 0.2470579220100149
$ws.Range("O20").Value = 4.202643824197139

$ws.Range("B21").Value = 0.6785281462777277
$ws.Range("C21").Value = 0.1607266961177345
$ws.Range("E21").Value = 0.2185050480002459
$ws.Range("F21").Value = 2.225558008341238
$ws.Range("G21").Value = 0.00248441066353773
$ws.Range("I21").Value = 1.109246630546295
$ws.Range("J21").Value = 0.09051468666661577
$ws.Range("K21").Value = 0.414134243561989
$ws.Range("L21").Value = 0.4604683861488184
$ws.Range("M21").Value = 0.2590341406579384
$ws.Range("O21").Value = 4.172443332422148

$ws.Range("B22").Value = 0.7114392231939632
$ws.Range("C22").Value = 0.1612191797216695
$ws.Range("E22").Value = 0.2190493580333772
$ws.Range("F22").Value = 2.227122764882054
$ws.Range("G22").Value = 0.002482790125455843
$ws.Range("I22").Value = 1.103352713715182
$ws.Range("J22").Value = 0.09208844227990909
$ws.Range("K22").Value = 0.4419609511279532
$ws.Range("L22").Value = 0.4662766087722474
$ws.Range("M22").Value = 0.2669204411730064
$ws.Range("O22").Value = 4.154416930572495

$ws.Range("B23").Value = 0.6938668261079783
$ws.Range("C23").Value = 0.1609561804653694
$ws.Range("E23").Value = 0.2187534936029252
$ws.Range("F23").Value = 2.226207526772228
$ws.Range("G23").Value = 0.002483649206333969
$ws.Range("I23").Value = 1.106460123345059
$ws.Range("J23").Value = 0.09124929399165183
$ws.Range("K23").Value = 0.4271116156824348
$ws.Range("L23").Value = 0.4631644773447761
$ws.Range("M23").Value = 0.2627060233568344
$ws.Range("O23").Value = 4.163881403493519

$ws.Range("B24").Value = 0.6274749212569191
$ws.Range("C24").Value = 0.1599635194008542
$ws.Range("E24").Value = 0.2177509059501119
$ws.Range("F24").Value = 2.224502404744854
$ws.Range("G24").Value = 0.002487031210400401
$ws.Range("I24").Value = 1.119070178585865
$ws.Range("J24").Value = 0.08805383951639811
$ws.Range("K24").Value = 0.3708244628175521
$ws.Range("L24").Value = 0.4516462375062531
$ws.Range("M24").Value = 0.2468631891812905
$ws.Range("O24").Value = 4.203166334382615

$ws.Range("B25").Value = 0.5562293587036606
$ws.Range("C25").Value = 0.1589005578844791
$ws.Range("E25").Value = 0.2169232774421026
$ws.Range("F25").Value = 2.226435450506614
$ws.Range("G25").Value = 0.002490955267055874
$ws.Range("I25").Value = 1.134470404294333
$ws.Range("J25").Value = 0.08457095722375385
$ws.Range("K25").Value = 0.3100264616405184
$ws.Range("L25").Value = 0.4398019839308063
$ws.Range("M25").Value = 0.2300342105665578
$ws.Range("O25").Value = 4.252892392011745
